{"js": "// Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" footer block\n// (and its leading blank paragraph, and the copyright line right after\n// it) that used to follow the \"LOQ4085: ... (Requisito fraco)\" paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst ANCHOR_TEXT = \"LOQ4085: Opera\u00e7\u00f5es Unit\u00e1rias I (Requisito fraco)\";\nconst JUPITER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_PREFIX = \"\u00a9 2020\";\n\n// Locate the anchor paragraph (\"LOQ4085: ...\") that precedes the block\n// we need to remove.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === ANCHOR_TEXT) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the anchor paragraph: \" + ANCHOR_TEXT);\n}\n\n// Collect the paragraphs to delete: starting right after the anchor,\n// an (optional) blank paragraph, followed by the \"Ver no Jupiter...\"\n// paragraph and the \"\u00a9 2020 ...\" paragraph. Stop as soon as we pass the\n// copyright paragraph so unrelated later content is left untouched.\nconst toRemove = [];\nlet sawJupiter = false;\nfor (let i = anchorIndex + 1; i < items.length; i++) {\n  const text = items[i].text;\n\n  if (text === JUPITER_TEXT) {\n    toRemove.push(items[i]);\n    sawJupiter = true;\n    continue;\n  }\n\n  if (text.indexOf(COPYRIGHT_PREFIX) === 0) {\n    toRemove.push(items[i]);\n    break;\n  }\n\n  if (!sawJupiter && text === \"\") {\n    // Leading blank paragraph before the \"Ver no Jupiter...\" line.\n    toRemove.push(items[i]);\n    continue;\n  }\n\n  // Hit something that isn't part of the block we're removing; stop.\n  break;\n}\n\n// Delete the collected paragraphs (reverse order is unnecessary since we\n// hold direct object references, but it's harmless and keeps things tidy).\nfor (let i = toRemove.length - 1; i >= 0; i--) {\n  toRemove[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter Salvar em pdf Salvar em docx\" footer block\n# (its leading blank paragraph and the copyright line right after it)\n# that used to follow the \"LOQ4085: ... (Requisito fraco)\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the \"LOQ4085: ... (Requisito fraco)\" paragraph that anchors the\n# block we need to remove.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -like \"LOQ4085*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph (LOQ4085...)\"\n}\n\n# Walk forward from the anchor collecting the indices of the paragraphs\n# that make up the footer block to delete: an optional blank paragraph,\n# the \"Ver no Jupiter...\" paragraph, and the \"(c) 2020 ...\" paragraph.\n# Stop as soon as something outside that block is seen so unrelated\n# later content (e.g. the trailing page-break paragraph) is untouched.\n$toDelete = New-Object System.Collections.ArrayList\n$sawJupiter = $false\nfor ($i = $anchorIndex + 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($t -like \"Ver no Jupiter*\") {\n        [void]$toDelete.Add($i)\n        $sawJupiter = $true\n        continue\n    }\n\n    if ($t -like \"*2020*Contact*\") {\n        [void]$toDelete.Add($i)\n        break\n    }\n\n    if (-not $sawJupiter -and $t -eq \"\") {\n        [void]$toDelete.Add($i)\n        continue\n    }\n\n    break\n}\n\n# Delete from the highest index down so earlier indices stay valid as we\n# remove paragraphs.\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n    $idx = $toDelete[$j]\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
